# Insert a new row at position 78 (shifts existing rows 78..201 down to 79..202)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(78).Insert()

# Populate the newly inserted row 78 with the new data record
$ws.Range("A78").Value = 5
$ws.Range("B78").Value = "Macroferia Regional de Talca"
$ws.Range("C78").Value = "Maule"
$ws.Range("D78").Value = 44477
$ws.Range("E78").Value = 7
$ws.Range("F78").Value = 100112023
$ws.Range("G78").Value = "Brócoli"
$ws.Range("H78").Value = "Sin especificar"
$ws.Range("I78").Value = "Primera"
$ws.Range("J78").Value = 3000
$ws.Range("K78").Value = 600
$ws.Range("L78").Value = 600
$ws.Range("M78").Value = 600
$ws.Range("N78").Value = "$/unidad"
$ws.Range("O78").Value = "Región del Maule"
$ws.Range("P78").Value = 600
$ws.Range("Q78").Value = 1
$ws.Range("R78").Value = "Hortaliza"
